$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# Settings sheet: rename queue + add asset folder value
$wsSettings.Range("B2").Value = "AutoRecrut"
$wsSettings.Range("B3").Value = "bello"

# Constants sheet: add new email-related exception constants
# Names (column A) were entered first for all rows, then the values (column B).
$wsConstants.Range("A19").Value = "BException_Email_Body"
$wsConstants.Range("A20").Value = "BException_Email_Subject "
$wsConstants.Range("A21").Value = "SException_Email_Subject "
$wsConstants.Range("A22").Value = "SException_Email_Body"

$wsConstants.Range("B20").Value = "No Attachment Found"
$wsConstants.Range("B20").WrapText = $true

$wsConstants.Range("B19").Value = "Hello , Kindly note that the Email provided had no attachment added on it ,Kind Regards Admin"
$wsConstants.Range("B19").WrapText = $true

$wsConstants.Range("B21").Value = "Hi , Certain Errors were experienced In the system"

$wsConstants.Range("B22").Value = "System Error"
